$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''331.96'
$ws.Range('E2').Value = '''1.02%'
$ws.Range('D3').Value = '''45.89'
$ws.Range('E3').Value = '''4.36%'
$ws.Range('D4').Value = '''5.622'
$ws.Range('E4').Value = '''1.99%'
$ws.Range('D5').Value = '''0.08363'
$ws.Range('E5').Value = '''4.20%'
$ws.Range('D6').Value = '''2.054'
$ws.Range('E6').Value = '''3.20%'
$ws.Range('B7').Value = 'MXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D7').Value = '''0.9829'
$ws.Range('E7').Value = '''3.49%'
$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').Value = '''2.594'
$ws.Range('E8').Value = '''1.16%'
$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D9').Value = '''0.1157'
$ws.Range('E9').Value = '''1.53%'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = '''0.1926'
$ws.Range('E10').Value = '''3.01%'
$ws.Range('B11').Value = 'MCDex'
$ws.Range('C11').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D11').Value = '''10.37'
$ws.Range('E11').Value = '''-3.41%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '''0.09983'
$ws.Range('E12').Value = '''1.40%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '''0.04665'
$ws.Range('E13').Value = '''-0.77%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '''0.1058'
$ws.Range('E14').Value = '''-0.63%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '''0.001275'
$ws.Range('E15').Value = '''0.71%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '''0.006086'
$ws.Range('E16').Value = '''2.69%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '''3.374'
$ws.Range('E17').Value = '''0.49%'
$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').Value = '''4.480'
$ws.Range('E18').Value = '''3.28%'
$ws.Range('D19').Value = '''0.3365'
$ws.Range('E19').Value = '''-3.17%'
$ws.Range('E20').Value = '''-1.25%'
$ws.Range('E21').Value = '''4.23%'
$ws.Range('D22').Value = '''0.04206'
$ws.Range('E22').Value = '''3.25%'
$ws.Range('E23').Value = '''4.24%'
$ws.Range('D24').Value = '''0.004623'
$ws.Range('E24').Value = '''6.67%'
$ws.Range('D26').Value = '''0.0003748'
$ws.Range('E26').Value = '''0.10%'
$ws.Range('D38').Value = '''0.02778'
$ws.Range('E38').Value = '''7.37%'
$ws.Range('D39').Value = '''0.05794'
$ws.Range('E39').Value = '''2.63%'
$ws.Range('D40').Value = '''0.007764'
$ws.Range('E40').Value = '''3.82%'
$ws.Range('D41').Value = '''0.1433'
$ws.Range('E41').Value = '''2.58%'
$ws.Range('D42').Value = '''0.007281'
$ws.Range('E42').Value = '''-3.53%'
$ws.Range('D43').Value = '''0.002014'
$ws.Range('E43').Value = '''-0.01%'
$ws.Range('D44').Value = '''0.008100'
$ws.Range('E44').Value = '''-5.80%'
$ws.Range('D45').Value = '''0.3410'
$ws.Range('D46').Value = '''0.00007333'
$ws.Range('E46').Value = '''2.71%'
$ws.Range('E47').Value = '''0.22%'
$ws.Range('D48').Value = '''0.0005812'
$ws.Range('E48').Value = '''0.01%'
$ws.Range('B49').Value = 'BOLO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range('D49').Value = '''0.003499'
$ws.Range('E49').Value = '''-5.21%'
$ws.Range('B50').Value = 'CoinbaseStockToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range('D50').Value = '''0.003505'
$ws.Range('E50').Value = '''-0.68%'
$ws.Range('D51').Value = '''0.00002104'
$ws.Range('E51').Value = '''0.22%'

# Reset number format to default (no explicit style) for cells written
# with a leading apostrophe so the saved XML has no style index,
# matching the original plain inline-string cells.
$ws.Range('D2,E2,D3,E3,D4,E4,D5,E5,D6,E6,D7,E7,D8,E8,D9,E9,D10,E10,D11,E11,D12,E12,D13,E13,D14,E14,D15,E15,D16,E16,D17,E17,D18,E18,D19,E19,E20,E21,D22,E22,E23,D24,E24,D26,E26,D38,E38,D39,E39,D40,E40,D41,E41,D42,E42,D43,E43,D44,E44,D45,D46,E46,E47,D48,E48,D49,E49,D50,E50,D51,E51').Style = 'Normal'
